$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new title paragraph at the very top of the document.
# ---------------------------------------------------------------------
$firstPara = $d.Paragraphs(1)
$firstPara.Range.InsertParagraphBefore()
$titleRange = $d.Paragraphs(1).Range
$titleXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pBdr>
      <w:bottom w:val="double" w:sz="6" w:space="1" w:color="auto"/>
    </w:pBdr>
    <w:jc w:val="center"/>
    <w:rPr>
      <w:rFonts w:ascii="Aharoni" w:hAnsi="Aharoni" w:cs="Aharoni"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="44"/>
      <w:szCs w:val="44"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aharoni" w:hAnsi="Aharoni" w:cs="Aharoni"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="44"/>
      <w:szCs w:val="44"/>
    </w:rPr>
    <w:t xml:space="preserve">Clinical Structured History taking Inventory for Children aged </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aharoni" w:hAnsi="Aharoni" w:cs="Aharoni"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="60"/>
      <w:szCs w:val="60"/>
    </w:rPr>
    <w:t>4-8</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aharoni" w:hAnsi="Aharoni" w:cs="Aharoni"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="72"/>
      <w:szCs w:val="72"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aharoni" w:hAnsi="Aharoni" w:cs="Aharoni"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="44"/>
      <w:szCs w:val="44"/>
    </w:rPr>
    <w:t>years to explore Cerebral Visual Impairment</w:t>
  </w:r>
</w:p>
'@
$titleRange.InsertXML($titleXml)

# ---------------------------------------------------------------------
# 2) Append four new question/answer paragraphs after "Date of Birth"
#    (the current last paragraph in the document body).
# ---------------------------------------------------------------------
function Add-QAParagraph([string]$question, [string]$field) {
    $last = $d.Paragraphs($d.Paragraphs.Count)
    $last.Range.InsertParagraphAfter()
    $newRange = $d.Paragraphs($d.Paragraphs.Count).Range
    $xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>$question</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:br/>
  </w:r>
  <w:r>
    <w:t>{</w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>$field</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>}</w:t>
  </w:r>
</w:p>
"@
    $newRange.InsertXML($xml)
}

Add-QAParagraph "How many weeks into your pregnancy were you when you gave birth?" "pregWeeks"
Add-QAParagraph "Were there any problems at the time of birth? If so, please describe:" "birthProblems"
Add-QAParagraph "Has your child had any conditions affecting the eyes or brain? If so, please describe:" "affectingConditions"
Add-QAParagraph "Do you have any concerns about your child’s vision? If so, what are they?" "visionConcerns"

# ---------------------------------------------------------------------
# 3) Un-hide the "Default Paragraph Font" character style.
# ---------------------------------------------------------------------
$style = $d.Styles("Default Paragraph Font")
$style.Hidden = $false

Write-Host "done"
